$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Cells that become literal text placeholders ('0') ---
# Copy from C14, which already holds text '0' with the matching style,
# so the destination keeps the same style + shared-string text type.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("C29"))

# --- Cells that become literal text placeholders ('***.*') ---
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("E14").Copy($ws.Range("E26"))

# --- Plain numeric value updates ---
$ws.Range("L14").Value = -33.333333333333
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 75
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 87.5
$ws.Range("L16").Value = 127.272727272727
$ws.Range("M16").Value = 47.058823529411
$ws.Range("N16").Value = -82.837528604119
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -21.052631578947
$ws.Range("I17").Value = 112
$ws.Range("J17").Value = 115
$ws.Range("K17").Value = -2.608695652173
$ws.Range("M17").Value = 38.271604938271
$ws.Range("N17").Value = -34.502923976608
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -58.333333333333
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = -15.151515151515
$ws.Range("L18").Value = 68
$ws.Range("M18").Value = 21.739130434782
$ws.Range("N18").Value = -80.327868852459
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -7.142857142857
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -11.864406779661
$ws.Range("I19").Value = 382
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = 14.371257485029
$ws.Range("L19").Value = 98.958333333333
$ws.Range("M19").Value = 34.982332155477
$ws.Range("N19").Value = -55.684454756380
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 250
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 75
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = 211.111111111111
$ws.Range("N20").Value = -83.815028901734
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -18.75
$ws.Range("I21").Value = 686
$ws.Range("J21").Value = 616
$ws.Range("K21").Value = 11.363636363636
$ws.Range("L21").Value = 75.447570332480
$ws.Range("M21").Value = 38.306451612903
$ws.Range("N21").Value = -67.003367003367
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -40
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = -53.333333333333
$ws.Range("L22").Value = -22.222222222222
$ws.Range("M22").Value = -22.222222222222
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -15.384615384615
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 22.222222222222
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 28.571428571428
$ws.Range("G24").Value = 61
$ws.Range("H24").Value = 8.196721311475
$ws.Range("I24").Value = 569
$ws.Range("J24").Value = 435
$ws.Range("K24").Value = 30.804597701149
$ws.Range("L24").Value = 60.281690140845
$ws.Range("M24").Value = -9.538950715421
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 39.130434782608
$ws.Range("I25").Value = 196
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = -2
$ws.Range("L25").Value = 21.739130434782
$ws.Range("M25").Value = 28.104575163398
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 39
$ws.Range("J27").Value = 61
$ws.Range("K27").Value = -36.065573770491
$ws.Range("L27").Value = -26.415094339622
$ws.Range("G30").Value = 1
